# Applies the "fixed env file stuff" edit:
#  - Inserts three new columns (H,I,J) for "Metazoa_*" data right before the
#    existing "Vertebrata_*" columns, which shift from H,I,J to K,L,M.
#  - Updates the row-1 headers accordingly.
#  - Populates the new image_file hyperlink formula for rows 4 and 5 (col G).
#  - Refreshes the Vertebrata values for rows 4 and 5 (now in K/L/M) and fills
#    in the new Metazoa values for rows 4 and 5 (H/I/J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert three blank columns before the old column H (Vertebrata_*) ---
# Old layout:  H=Vertebrata_aln_property_entropy_z_score, I=Vertebrata_aln_slice_view, J=Vertebrata_cons_string
# New layout:  H,I,J = Metazoa_*  and K,L,M = Vertebrata_* (shifted right by 3)
$ws.Range("H1:J1").Insert(-4161)

# --- Row 1 headers ---
$ws.Range("H1").Value = "Metazoa_aln_asym_sum_of_pairs_z_score"
$ws.Range("I1").Value = "Metazoa_aln_slice_view"
$ws.Range("J1").Value = "Metazoa_cons_string"
$ws.Range("K1").Value = "Vertebrata_aln_asym_sum_of_pairs_z_score"

# --- Row 4 (reference_index = 2) ---
$ws.Range("G4").Formula = '=HYPERLINK("/Users/jackson/Dropbox (MIT)/work/07-SLiM_bioinformatics/05-conservation_pipeline/examples/table_annotation/conservation_analysis/2-9606_0_002f40/2-9606_0002f40-aln_asym_sum_of_pairs_og_level_score_screen.png")'
$ws.Range("H4").Value = -0.2623404708668192
$ws.Range("I4").Formula = '=HYPERLINK("conservation_analysis/2-9606_0_002f40/2-9606_0002f40-Metazoa_aln_slice.html")'
$ws.Range("J4").Value = "__SP_P_____"
$ws.Range("K4").Value = -0.2623404708668192
$ws.Range("M4").Value = "__SP_P_____"

# --- Row 5 (reference_index = 3) ---
$ws.Range("G5").Formula = '=HYPERLINK("/Users/jackson/Dropbox (MIT)/work/07-SLiM_bioinformatics/05-conservation_pipeline/examples/table_annotation/conservation_analysis/3-9606_0_002f40/3-9606_0002f40-aln_asym_sum_of_pairs_og_level_score_screen.png")'
$ws.Range("H5").Value = -1.164739051311995
$ws.Range("I5").Formula = '=HYPERLINK("conservation_analysis/3-9606_0_002f40/3-9606_0002f40-Metazoa_aln_slice.html")'
$ws.Range("J5").Value = "_______EE___"
$ws.Range("K5").Value = -1.164739051311995
$ws.Range("M5").Value = "_______EE___"

Write-Host "Applied Metazoa/Vertebrata column split edit"
